$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 4055.1072
$ws.Range("I19").Value = 6244.9414
$ws.Range("J19").Value = 670.8182
$ws.Range("K19").Value = 6244.9414
$ws.Range("L19").Value = 670.8182
$ws.Range("M19").Value = -6069.9414
$ws.Range("N19").Value = -1020.8182
# Row 40
$ws.Range("H40").Value = 4258.4
$ws.Range("I40").Value = 800
$ws.Range("J40").Value = 4790.4614
$ws.Range("K40").Value = 800
$ws.Range("L40").Value = 4790.4614
$ws.Range("M40").Value = -625
$ws.Range("N40").Value = -5140.4614
# Row 138
$ws.Range("H138").Value = 2722.9868
$ws.Range("I138").Value = 1957.2632
$ws.Range("J138").Value = 2978.228
$ws.Range("K138").Value = 5871.7896
$ws.Range("L138").Value = 8934.684000000001
$ws.Range("M138").Value = -731.7896000000001
$ws.Range("N138").Value = -19214.684

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 1858.8334
$ws.Range("I110").Value = 2135
$ws.Range("J110").Value = 1306.5
$ws.Range("K110").Value = 2135
$ws.Range("L110").Value = 1306.5
$ws.Range("M110").Value = -90
$ws.Range("N110").Value = -5396.5
# Row 132
$ws.Range("H132").Value = 2942.3052
$ws.Range("I132").Value = 2782
$ws.Range("J132").Value = 3338.353
$ws.Range("K132").Value = 8346
$ws.Range("L132").Value = 10015.059
$ws.Range("M132").Value = -5816
$ws.Range("N132").Value = -15075.059

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2597.6482
$ws.Range("I105").Value = 2333.2273
$ws.Range("J105").Value = 3761.1
$ws.Range("K105").Value = 2333.2273
$ws.Range("L105").Value = 3761.1
$ws.Range("M105").Value = -586.2273
$ws.Range("N105").Value = -7255.1

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2512.4458
$ws.Range("I31").Value = 1863.6038
$ws.Range("J31").Value = 3658.7334
$ws.Range("K31").Value = 1863.6038
$ws.Range("L31").Value = 3658.7334
$ws.Range("M31").Value = -1568.6038
$ws.Range("N31").Value = -4248.7334
# Row 34
$ws.Range("H34").Value = 2512.4458
$ws.Range("I34").Value = 1863.6038
$ws.Range("J34").Value = 3658.7334
$ws.Range("K34").Value = 1863.6038
$ws.Range("L34").Value = 3658.7334
$ws.Range("M34").Value = -1661.6038
$ws.Range("N34").Value = -4062.7334
# Row 50
$ws.Range("H50").Value = 16000
$ws.Range("J50").Value = 16000
$ws.Range("L50").Value = 16000
$ws.Range("N50").Value = -17250
# Row 51
$ws.Range("H51").Value = 24875
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 24875
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 24875
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -26347
# Row 59
$ws.Range("H59").Value = 62450
$ws.Range("I59").Value = 40000
$ws.Range("J59").Value = 69933.336
$ws.Range("K59").Value = 40000
$ws.Range("L59").Value = 69933.336
$ws.Range("M59").Value = -38855
$ws.Range("N59").Value = -72223.336
# Row 60
$ws.Range("H60").Value = 29800
$ws.Range("J60").Value = 29800
$ws.Range("L60").Value = 29800
$ws.Range("N60").Value = -30822
# Row 61
$ws.Range("H61").Value = 24875
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 24875
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 24875
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -25571
# Row 68
$ws.Range("H68").Value = 29246.6
$ws.Range("I68").Value = 15000
$ws.Range("J68").Value = 32808.25
$ws.Range("K68").Value = 15000
$ws.Range("L68").Value = 32808.25
$ws.Range("M68").Value = -14251
$ws.Range("N68").Value = -34306.25
# Row 71
$ws.Range("H71").Value = 29246.6
$ws.Range("I71").Value = 15000
$ws.Range("J71").Value = 32808.25
$ws.Range("K71").Value = 45000
$ws.Range("L71").Value = 98424.75
$ws.Range("M71").Value = -41256
$ws.Range("N71").Value = -105912.75
# Row 74
$ws.Range("H74").Value = 20485.5
$ws.Range("J74").Value = 20485.5
$ws.Range("L74").Value = 20485.5
$ws.Range("N74").Value = -22233.5
# Row 77
$ws.Range("H77").Value = 20485.5
$ws.Range("J77").Value = 20485.5
$ws.Range("L77").Value = 61456.5
$ws.Range("N77").Value = -70192.5
# Row 132
$ws.Range("H132").Value = 16669584
$ws.Range("I132").Value = 62501268
$ws.Range("J132").Value = 3516.4546
$ws.Range("K132").Value = 187503804
$ws.Range("L132").Value = 10549.3638
$ws.Range("M132").Value = -187501274
$ws.Range("N132").Value = -15609.3638
# Row 141
$ws.Range("H141").Value = 25067.143
$ws.Range("J141").Value = 23918.46
$ws.Range("L141").Value = 23918.46
$ws.Range("N141").Value = -34278.46

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 655.3871
$ws.Range("I5").Value = 437.36365
$ws.Range("J5").Value = 1188.3334
$ws.Range("K5").Value = 1312.09095
$ws.Range("L5").Value = 3565.0002
$ws.Range("M5").Value = -1200.09095
$ws.Range("N5").Value = -3789.0002
# Row 38
$ws.Range("H38").Value = 48.11111
$ws.Range("J38").Value = 55
$ws.Range("L38").Value = 165
$ws.Range("N38").Value = -859
# Row 122
$ws.Range("H122").Value = 833.1667
$ws.Range("I122").Value = 427.33334
$ws.Range("J122").Value = 1239
$ws.Range("K122").Value = 3846.00006
$ws.Range("L122").Value = 11151
$ws.Range("M122").Value = -1396.00006
$ws.Range("N122").Value = -16051
# Row 132
$ws.Range("H132").Value = 5723.8096
$ws.Range("I132").Value = 5233.3335
$ws.Range("J132").Value = 5920
$ws.Range("K132").Value = 47100.0015
$ws.Range("L132").Value = 53280
$ws.Range("M132").Value = -44570.0015
$ws.Range("N132").Value = -58340
# Row 135
$ws.Range("H135").Value = 655.3871
$ws.Range("I135").Value = 437.36365
$ws.Range("J135").Value = 1188.3334
$ws.Range("K135").Value = 3936.27285
$ws.Range("L135").Value = 10695.0006
$ws.Range("M135").Value = -1401.27285
$ws.Range("N135").Value = -15765.0006
# Row 139
$ws.Range("H139").Value = 19233554
$ws.Range("I139").Value = 29412436
$ws.Range("J139").Value = 6777.778
$ws.Range("K139").Value = 88237308
$ws.Range("L139").Value = 20333.334
$ws.Range("M139").Value = -88232168
$ws.Range("N139").Value = -30613.334

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 6150
$ws.Range("I7").Value = 9800
$ws.Range("J7").Value = 2500
$ws.Range("K7").Value = 9800
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = -9688
$ws.Range("N7").Value = -2724
# Row 8
$ws.Range("H8").Value = 6150
$ws.Range("I8").Value = 9800
$ws.Range("J8").Value = 2500
$ws.Range("K8").Value = 9800
$ws.Range("L8").Value = 2500
$ws.Range("M8").Value = -9661
$ws.Range("N8").Value = -2778
# Row 107
$ws.Range("H107").Value = 8598.666999999999
$ws.Range("I107").Value = 10210.5
$ws.Range("K107").Value = 10210.5
$ws.Range("M107").Value = -8290.5
# Row 113
$ws.Range("H113").Value = 5108.0415
$ws.Range("I113").Value = 6989.6875
$ws.Range("J113").Value = 1344.75
$ws.Range("K113").Value = 6989.6875
$ws.Range("L113").Value = 1344.75
$ws.Range("M113").Value = -4819.6875
$ws.Range("N113").Value = -5684.75
# Row 132
$ws.Range("H132").Value = 4885.355
$ws.Range("I132").Value = 6166.7144
$ws.Range("J132").Value = 3830.1177
$ws.Range("K132").Value = 18500.1432
$ws.Range("L132").Value = 11490.3531
$ws.Range("M132").Value = -15970.1432
$ws.Range("N132").Value = -16550.3531

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 927.2727
$ws.Range("I46").Value = 925
$ws.Range("J46").Value = 928.5714
$ws.Range("K46").Value = 925
$ws.Range("L46").Value = 928.5714
$ws.Range("M46").Value = -737
$ws.Range("N46").Value = -1304.5714

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 5589
$ws.Range("J41").Value = 5589
$ws.Range("L41").Value = 5589
$ws.Range("N41").Value = -6369
# Row 45
$ws.Range("H45").Value = 5000
$ws.Range("J45").Value = 5000
$ws.Range("L45").Value = 5000
$ws.Range("N45").Value = -5982
# Row 74
$ws.Range("H74").Value = 7333.2
$ws.Range("J74").Value = 7333.2
$ws.Range("L74").Value = 7333.2
$ws.Range("N74").Value = -9205.200000000001
# Row 77
$ws.Range("H77").Value = 7333.2
$ws.Range("J77").Value = 7333.2
$ws.Range("L77").Value = 21999.6
$ws.Range("N77").Value = -31359.6
# Row 122
$ws.Range("H122").Value = 40333.617
$ws.Range("I122").Value = 54646
$ws.Range("J122").Value = 1485.7142
$ws.Range("K122").Value = 163938
$ws.Range("L122").Value = 4457.142599999999
$ws.Range("M122").Value = -161488
$ws.Range("N122").Value = -9357.142599999999

Write-Host "Applied all Unicorn_Profits updates."